$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 211; this shifts the existing rows 211:313 down to
# 212:314 and grows the used range / dimension to A1:R314 automatically.
$ws.Rows("211:211").Insert()

# Populate the newly inserted row 211 with the new data record.
$ws.Range("A211").Value = 10
$ws.Range("B211").Value = 'Vega Modelo de Temuco'
$ws.Range("C211").Value = 'La Araucanía'
$ws.Range("D211").Value = 44726
$ws.Range("E211").Value = 9
$ws.Range("F211").Value = 100112017
$ws.Range("G211").Value = 'Apio'
$ws.Range("H211").Value = 'Americana (o)'
$ws.Range("I211").Value = 'Primera'
$ws.Range("J211").Value = 60
$ws.Range("K211").Value = 10000
$ws.Range("L211").Value = 10000
$ws.Range("M211").Value = 10000
$ws.Range("N211").Value = '$/docena de matas'
$ws.Range("O211").Value = 'Provincia del Elquí'
$ws.Range("P211").Value = 1667
$ws.Range("Q211").Value = 6
$ws.Range("R211").Value = 'Hortaliza'
